# Slide 3 ("...") - shape id 6 "Rectangle 5" holds a Java code sample.
# Two small in-place text edits to existing runs:
#   1. "int" -> "double"            (declared type of `width`)
#   2. "%s" -> "%d" inside the printf format string
#
# The shape auto-fits its height to the text (a:spAutoFit). Re-measuring
# text through this host's text-metrics approximation does not reproduce
# PowerPoint's original layout exactly, so after the edits we restore the
# shape's AutoSize behaviour and exact height to avoid spurious size drift.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(4)

if ($sh.Name -ne "Rectangle 5") {
    throw "Unexpected shape: $($sh.Name)"
}

$tf = $sh.TextFrame
$tr = $tf.TextRange

$origAutoSize = $tf.AutoSize
$origHeight = $sh.Height

# Prevent the shape from auto-resizing while we edit the text.
$tf.AutoSize = 0

# --- Edit 1: "int" -> "double" ------------------------------------------
$full = $tr.Text
$idx = $full.IndexOf("int")
if ($idx -lt 0) { throw "Could not find 'int' run" }
$run = $tr.Characters($idx + 1, 3)
if ($run.Text -ne "int") { throw "Unexpected text at int location: $($run.Text)" }
$run.Text = "double"

# --- Edit 2: printf("Width: %s\n", width); -> ...%d... -------------------
$full = $tr.Text
$pidx = $full.IndexOf("printf")
if ($pidx -lt 0) { throw "Could not find 'printf' run" }
$start = $pidx + "printf".Length
$semi = $full.IndexOf(";", $start)
if ($semi -lt 0) { throw "Could not find terminating ';'" }
$len = $semi - $start + 1
$run2 = $tr.Characters($start + 1, $len)
# NB: TextRange.Text normalizes the stored curly quotes to straight quotes
# when read back, so compare/guard against the straight-quote form, but
# write the replacement using the original curly quotes (U+201C/U+201D) so
# the saved XML matches the source formatting exactly.
$expectedRead = "(`"Width: %s\n`", width);"
if ($run2.Text -ne $expectedRead) { throw "Unexpected printf run text: $($run2.Text)" }
$run2.Text = "(" + [char]8220 + "Width: %d\n" + [char]8221 + ", width);"

# Restore auto-fit behaviour and the shape's original height precisely.
$tf.AutoSize = $origAutoSize
$sh.Height = $origHeight
# Height is stored as a 32-bit float internally and gets truncated (not
# rounded) back to EMUs on save, so re-assigning the exact read-back value
# can land one EMU short of the original. Nudge it by a hair so the
# round-trip lands back on the original EMU value.
$sh.Height = $origHeight + 0.00004
